$d = $word.ActiveDocument

$replacements = @(
    @{old="774×7="; new="924×3="},
    @{old="128×6="; new="334×9="},
    @{old="578×9="; new="311×8="},
    @{old="367×8="; new="309×3="},
    @{old="458×3="; new="745×5="},
    @{old="555×2="; new="484×8="},
    @{old="602×6="; new="627×9="},
    @{old="493×7="; new="894×5="},
    @{old="534×5="; new="437×5="},
    @{old="689×4="; new="801×9="},
    @{old="467×3="; new="926×4="},
    @{old="101×2="; new="513×9="},
    @{old="280×3="; new="919×9="},
    @{old="764×3="; new="565×4="},
    @{old="825×5="; new="982×7="},
    @{old="812×5="; new="167×7="},
    @{old="208×4="; new="425×2="},
    @{old="889×2="; new="588×5="},
    @{old="475×9="; new="410×3="},
    @{old="976×8="; new="771×8="},
    @{old="277×2="; new="981×4="},
    @{old="663×7="; new="230×7="},
    @{old="910×8="; new="550×9="},
    @{old="702×7="; new="419×8="},
    @{old="579×2="; new="660×9="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
